# Add a new bullet (list item) right after the "simeonvisser" hyperlink
# paragraph, containing a hyperlink to the Selenium Grid "run_the_demo.html"
# page in the git history, matching the sibling bullets' list formatting.

$d = $word.ActiveDocument

$url = "http://grid.selenium.googlecode.com/git-history/5b3078f83c748db7f9061a39b908daf877c7c55b/src/main/webapp/run_the_demo.html"

# Locate the paragraph whose text is the simeonvisser hyperlink (the last
# bullet in that list) by searching for its URL text.
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*simeonvisser.com*") {
        $anchorPara = $p
    }
}

# Insert a brand-new paragraph right after it, inheriting the same
# (ListParagraph / numId 2) list formatting.
$endRange = $anchorPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# The paragraph that was just created.
$newParaIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $newPara.Range

# Put the URL text in place, then trim back to an insertion point so the
# hyperlink can be created fresh (avoids leaving a stray empty run behind).
$newRange.Text = $url
$textRange = $d.Range($newRange.Start, $newRange.Start + $url.Length)
$textRange.Delete()

$insertionPoint = $d.Range($textRange.Start, $textRange.Start)
$d.Hyperlinks.Add($insertionPoint, $url) | Out-Null
